$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at N (shifts old N "目前利率" -> O, old O "調整後利率" -> P)
$ws.Columns("N:N").Insert()

# New header text for the inserted column
$ws.Range("N1").Value = "利率種類"

# Column widths: M stays narrower, new N column gets its own width,
# old widths for O/P (shifted from former N/O) are already carried over by Insert.
$ws.Columns("M:M").ColumnWidth = 10
$ws.Columns("N:N").ColumnWidth = 9.6

# Update the hidden _FilterDatabase defined name range to cover the new column
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=正常件!`$A`$1:`$P`$1"
    }
}

# Restore reference style (workbook calcPr no longer flagged R1C1)
$excel.ReferenceStyle = 1

# Leave the same cell selected as after the author's edit
$ws.Range("M7").Select()
